# TODO: new sample row appended to the "Account Information" sheet
# (placeholder "a" values plus a double-to-long-converted numeric id).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Account Information")

$ws1.Range("A3").Value = "a"
$ws1.Range("B3").Value = "a"
$ws1.Range("C3").Value = "a"
$ws1.Range("D3").Value = "a"
$ws1.Range("E3").Value = "a"
$ws1.Range("F3").Value = 89128184
